# US 3.3 commit files
# - Update the "Discount Rate" label on the VBDR sheet to clarify units.
# - Restore the cell selection (B2) that was recorded on the VBDR sheet view
#   the last time it was active, while leaving the "About" sheet as the
#   active/selected tab in the workbook (matching the saved file's state).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsVBDR  = $wb.Worksheets.Item("VBDR")

# Clarify that the discount rate is a dimensionless figure.
$wsVBDR.Range("B1").Value = "Discount Rate (dimensionless)"

# Record the selection on the VBDR sheet (cell B2), then switch back to the
# "About" sheet so it remains the active tab in the saved workbook.
$wsVBDR.Range("B2").Select()
$wsAbout.Activate()
